$wb = $excel.ActiveWorkbook

# --- Sheet "建物" (building): property_category column (I) was wrongly "land" ---
# Fix rows 2-7, column I, to read "building" instead of "land"
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 7; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value2 = "building"
}

# --- Sheet "汽車" (car): property_category column (H) was wrongly "land" ---
# Fix row 2, column H, to read "car" instead of "land"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value2 = "car"
